# "Generate Report for Handback"
# The CI tool re-generated the localization-status report after a handback:
#   - Status text for the two rows moved from "In Translation" to
#     "Handed back: in sync with en-US" on both language sheets.
#   - Each language sheet now also records the "Latest Target File" (a
#     hyperlinked .md name) and the "Latest Handback File" (.xlf name) for
#     both rows, plus a refreshed "Latest Handback DateTime".
#   - A few report columns were widened to comfortably fit the long file
#     names / guids that show up in them.

$wb = $excel.ActiveWorkbook

$mdUrlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09d039794719bed7ccf4cb5e0690261426f151a0/e2e/96f4894f-55c1-4e82-8f3d-3159110976b9.md"
$mdUrlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09d039794719bed7ccf4cb5e0690261426f151a0/e2e/b17468f1-c27e-4cc3-a404-f7ff7c40df27.md"
$mdNameA = "96f4894f-55c1-4e82-8f3d-3159110976b9.md"
$mdNameB = "b17468f1-c27e-4cc3-a404-f7ff7c40df27.md"

$newStatus = "Handed back: in sync with en-US"

function Set-HandbackLink {
    param($ws, [string]$cellRef, [string]$url, [string]$display)

    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $display) | Out-Null
    # Match the workbook's existing custom hyperlink look (underline + the
    # custom blue used by the A2/A3 hyperlinks) instead of the theme default
    # that a fresh "Hyperlink" style would apply.
    $ws.Range($cellRef).Font.Underline = 2
    $ws.Range($cellRef).Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Overview sheet: just a couple of columns growing wider in the refresh.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.08
$wsOverview.Columns.Item(6).ColumnWidth = 29.08

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(3).ColumnWidth = 29.08
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

Set-HandbackLink $wsZh "I2" $mdUrlA $mdNameA
Set-HandbackLink $wsZh "I3" $mdUrlB $mdNameB

$wsZh.Range("J2").Value = "96f4894f-55c1-4e82-8f3d-3159110976b9.416bc9c0f1592eddcfb7a6a058e633b58516e098.zh-cn.xlf"
$wsZh.Range("J3").Value = "b17468f1-c27e-4cc3-a404-f7ff7c40df27.46ae26565d58529d6750f5da8b18faa1b1e51ed3.zh-cn.xlf"

$wsZh.Range("K2").Value = "2016-08-23 00:23:38"
$wsZh.Range("K3").Value = "2016-08-23 00:23:38"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(3).ColumnWidth = 29.08
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

Set-HandbackLink $wsDe "I2" $mdUrlA $mdNameA
Set-HandbackLink $wsDe "I3" $mdUrlB $mdNameB

$wsDe.Range("J2").Value = "96f4894f-55c1-4e82-8f3d-3159110976b9.416bc9c0f1592eddcfb7a6a058e633b58516e098.de-de.xlf"
$wsDe.Range("J3").Value = "b17468f1-c27e-4cc3-a404-f7ff7c40df27.46ae26565d58529d6750f5da8b18faa1b1e51ed3.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-23 00:23:45"
$wsDe.Range("K3").Value = "2016-08-23 00:23:45"
